$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: "Correct_answer" header + per-row answer key ("l" for the
# Purple rows 2-7, "s" for the Blue rows 8-101).
$ws.Range("D1").Value = "Correct_answer"
$ws.Range("D2:D7").Value = "l"
$ws.Range("D8:D101").Value = "s"

# Match the author's final selection/active cell on the sheet.
$ws.Range("D8:D101").Select()
